# Generate Report for Handback
#
# This marks the d3bd625d-22f1-4cb9-9cd4-7fb766a6c0f8.md file as handed back
# (in sync with en-US) for both the zh-cn and de-de localization sheets, and
# records the handback target/handback-file/handback-datetime for each.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheetInfo = @(
    @{ Name = "zh-cn"; HandbackDateTime = "2016-02-18 09:18:36" },
    @{ Name = "de-de"; HandbackDateTime = "2016-02-18 09:18:58" }
)

# --- Overview sheet: update the status shown for the handed-off file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# --- Per-locale detail sheets ---
foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Look up the existing hyperlinks on A2 (source .md file) and C2 (handoff .xlf)
    # so the new "Latest Target File" (E2) and "Latest Handback File" (F2) columns
    # can point at the same targets/text.
    $sourceAddress = ""
    $sourceDisplay = ""
    $handoffAddress = ""
    $handoffDisplay = ""
    foreach ($hl in $ws.Hyperlinks) {
        $refAddr = $hl.Range.Address()
        if ($refAddr -eq "`$A`$2") {
            $sourceAddress = $hl.Address()
            $sourceDisplay = $hl.TextToDisplay()
        }
        if ($refAddr -eq "`$C`$2") {
            $handoffAddress = $hl.Address()
            $handoffDisplay = $hl.TextToDisplay()
        }
    }

    # Status -> Handed back
    $ws.Range("B2").Value = $newStatus

    # Latest Target File (E2) / Latest Handback File (F2): same file references
    # as the source (A2) and handoff (C2) columns respectively.
    $ws.Range("E2").Value = $sourceDisplay
    $ws.Range("F2").Value = $handoffDisplay
    $ws.Hyperlinks.Add($ws.Range("E2"), $sourceAddress, "", "", $sourceDisplay)
    $ws.Hyperlinks.Add($ws.Range("F2"), $handoffAddress, "", "", $handoffDisplay)

    # Match the existing hyperlink look (underlined, cornflower blue) used by
    # the other hyperlinked cells (A2/C2/A3) instead of the default theme
    # hyperlink style that Hyperlinks.Add applies.
    foreach ($col in @("E2", "F2")) {
        $fnt = $ws.Range($col).Font
        $fnt.Name = "Calibri"
        $fnt.Size = 11
        $fnt.Underline = 2
        $fnt.Color = 15570276
    }

    # Latest Handback DateTime (G2): now populated with the handback time.
    $ws.Range("G2").Value = $info.HandbackDateTime
}
